# Applies the "experimental update - few major features have been implemented"
# change set to the ingredients_matches workbook.
#
# Summary of the edit (derived from the canonical OOXML diff):
#  - sheetView: scroll position / selection moved from D1/L8 to C1/K7
#  - a handful of new "MATCHES" entries are written into the grid (new
#    shared strings), and a few previously-blank "no match" (red-filled)
#    placeholder cells gain new blank siblings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- view / selection -------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 3   # "C1" is column 3
$ws.Range("K7").Select() | Out-Null

# --- cell value edits ---------------------------------------------------
# New ingredient names are entered in the same order the author typed them
# so newly-created shared-string entries land at the same indices as the
# canonical workbook.
$ws.Range("N2").Value = "sweet_dew"
$ws.Range("L9").Value = "yam_bread"
$ws.Range("I5").Value = "frost_pie"
$ws.Range("G6").Value = "berry_icecream"
$ws.Range("H6").Value = "berry_jam"
$ws.Range("I6").Value = "nutty_cake"
$ws.Range("J6").Value = "berry_cake"
$ws.Range("I8").Value = "honey_pancakes"
$ws.Range("F5").ClearFormats()
$ws.Range("F5").Value = "berry_juice"

$ws.Range("J8").Value = "honey_cookies"
$ws.Range("M9").Value = "frost_pie"
$ws.Range("N9").Value = "cake_mix"
$ws.Range("O9").Value = "cake_mix"

# J10 previously held an empty "no match" (red-filled) placeholder; it now
# carries a real value and loses that formatting entirely.
$ws.Range("J10").ClearFormats()
$ws.Range("J10").Value = "cake"

# --- new blank "no match" placeholder cells (red fill, like the others) -
$redFill = $ws.Range("E2").Interior.Color
$ws.Range("G4").Interior.Color = $redFill
$ws.Range("G5").Interior.Color = $redFill
$ws.Range("J9").Interior.Color = $redFill
